$wb = $excel.ActiveWorkbook

# --- Update header on "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Update header on "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$data = @(
    @(44941.99999999999, 151, 85.87733251249436, 211.2870138597669),
    @(44955.99999999999, 99, 37.56765679643895, 159.3481385802377),
    @(44969.99999999999, 48, -15.89956262298243, 112.4856180525011),
    @(44976.99999999999, 22, -39.63048104324997, 86.74915544859414),
    @(44983.99999999999, 0, -68.89068895194825, 59.79850961480623),
    @(44990.99999999999, 0, -89.52397768419854, 38.0245509798093),
    @(44997.99999999999, 0, -120.2806341749897, 7.859989729836815),
    @(45004.99999999999, 0, -141.8798816626605, -14.03276839020511),
    @(45011.99999999999, 0, -173.3136389289872, -43.72948313026825),
    @(45018.99999999999, 0, -196.3814603900915, -68.63839941727625),
    @(45025.99999999999, 0, -219.4486659006959, -93.15976452261846),
    @(45032.99999999999, 0, -251.6858348271946, -120.6179434138514)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Range("A$row").Value = $r[0]
    $wsForecast.Range("A$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Range("B$row").Value = $r[1]
    $wsForecast.Range("C$row").Value = $r[2]
    $wsForecast.Range("D$row").Value = $r[3]
    $row++
}
